# Generate Report for Handback
# Refresh the timestamps recorded for the 908ad72e-1216-4cf0-8b8f-d81bd173824d
# handback row across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the 908ad72e... row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-28 04:44:58"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 908ad72e... row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-28 04:44:53"
$wsZhCn.Range("K4").Value = "2016-08-28 04:45:16"

# de-de sheet: "Correspond Handoff Datetime" (shares the same value as the
# Overview sheet's G4) / "Correspond Handback DateTime" for the
# 908ad72e... row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-28 04:44:58"
$wsDeDe.Range("K4").Value = "2016-08-28 04:45:22"
